$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.554.47'
$ws.Range('E2').Value = '  +3.49%  '
$ws.Range('D3').Value = '1.830.34'
$ws.Range('E3').Value = '  +5.00%  '
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').Value = "'344.53"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.09%  '
$ws.Range('D6').Value = "'0.9997"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.08%  '
$ws.Range('D7').Value = "'0.3836"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.26%  '
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').Value = "'0.3555"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.47%  '
$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').Value = "'50.22"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.63%  '
$ws.Range('D10').Value = "'1.246"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +4.24%  '
$ws.Range('D11').Value = "'0.07779"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +4.11%  '
$ws.Range('E12').Value = '  -0.07%  '
$ws.Range('D13').Value = "'22.28"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +8.90%  '
$ws.Range('D14').Value = "'6.629"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.44%  '
$ws.Range('D15').Value = '1.830.29'
$ws.Range('E15').Value = '  +5.18%  '
$ws.Range('E16').Value = '  +1.82%  '
$ws.Range('D17').Value = "'0.00001129"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.78%  '
$ws.Range('D18').Value = "'0.06753"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.81%  '
$ws.Range('D19').Value = "'87.03"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.06%  '
$ws.Range('D20').Value = "'1.000"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('D21').Value = "'17.64"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.22%  '
$ws.Range('D22').Value = "'6.580"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +6.44%  '
$ws.Range('E23').Value = '  +0.79%  '
$ws.Range('D24').Value = '27.567.39'
$ws.Range('E24').Value = '  +3.56%  '
$ws.Range('D25').Value = "'2.469"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.33%  '
$ws.Range('D26').Value = "'2.731"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +10.30%  '
$ws.Range('D27').Value = "'22.16"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +12.67%  '
$ws.Range('D28').Value = "'1.504"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +6.20%  '
$ws.Range('D29').Value = "'153.74"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.29%  '
$ws.Range('D30').Value = '2.036.49'
$ws.Range('E30').Value = '  +5.42%  '
$ws.Range('D31').Value = "'135.67"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.81%  '
$ws.Range('D32').Value = "'6.391"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +4.71%  '
$ws.Range('D33').Value = "'4.094"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.99%  '
$ws.Range('D34').Value = "'14.02"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +7.54%  '
$ws.Range('D35').Value = "'0.08825"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.07%  '
$ws.Range('D36').Value = "'1.695"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.72%  '
$ws.Range('D37').Value = "'5.665"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +4.32%  '
$ws.Range('D38').Value = "'0.7102"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +13.54%  '
$ws.Range('D39').Value = "'9.125"
$ws.Range('D39').Style = 'Normal'
$ws.Range('D40').Value = "'0.06542"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.68%  '
$ws.Range('B41').Value = 'Algorand'
$ws.Range('C41').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D41').Value = "'0.2264"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.74%  '
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').Value = "'0.02412"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +2.15%  '
$ws.Range('D43').Value = "'1.319"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +7.34%  '
$ws.Range('D44').Value = "'14.81"
$ws.Range('D44').Style = 'Normal'
$ws.Range('D45').Value = "'0.6654"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +9.77%  '
$ws.Range('D46').Value = "'0.9996"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.07%  '
$ws.Range('D47').Value = "'3.962"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.80%  '
$ws.Range('D48').Value = "'2.199"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +6.65%  '
$ws.Range('D49').Value = "'133.50"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.26%  '
$ws.Range('E50').Value = '  +0.45%  '
$ws.Range('D51').Value = "'81.42"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.58%  '
